$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP column (F) for rows 2-6 so that all server entries point at
# 127.0.0.1 (previously 192.168.1.113-192.168.1.117). Row 2 is included so
# that the old distinct IP strings collapse into a single shared string.
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F3").Value = "127.0.0.1"
$ws.Range("F4").Value = "127.0.0.1"
$ws.Range("F5").Value = "127.0.0.1"
$ws.Range("F6").Value = "127.0.0.1"

# Move the active cell selection to F14, as in the edited workbook.
$ws.Range("F14").Select()
